# Insert a new caption row above the old row 2 (the thin divider row) that
# labels the table's units -- "(in percent)" -- in the sheet's usual three
# languages (Kyrgyz / Russian / English). Every row from the old row 2
# downward shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 2; existing rows 2.. shift to 3..
$ws.Rows.Item(2).Insert() | Out-Null

# Populate the new row 2 with the localized "(in percent)" captions.
$ws.Cells.Item(2, 1).Value = "(пайыз менен)"
$ws.Cells.Item(2, 2).Value = "(в процентах)"
$ws.Cells.Item(2, 3).Value = "(in percent)"

# The insert copied row 1's formatting into D2:E2 even though they stay
# empty -- drop them so the row only carries the three caption cells.
$ws.Range("D2:E2").Clear() | Out-Null

# Build the exact target look (italic 8pt Times New Roman, centered) on an
# out-of-the-way scratch cell first, then copy/paste that formatting onto
# the caption cells in one shot -- this reuses a single shared style entry
# for A2:C2 instead of stamping a separate one per cell.
$scratch = $ws.Cells.Item(100, 100)
$scratch.Font.Italic = $true
$scratch.Font.Size = 8
$scratch.Font.Name = "Times New Roman"
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108
$scratch.WrapText = $false

$scratch.Copy() | Out-Null
$ws.Range("A2:C2").PasteSpecial(-4122) | Out-Null
$scratch.Clear() | Out-Null

$excel.CutCopyMode = $false

# Move the selection cursor to match the saved file.
$ws.Range("B10").Select() | Out-Null
